$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.02%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.76%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.570"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.14%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08054"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.57%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.971"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.47%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.320"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.07%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'-8.50%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9452"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.75%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1168"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.47%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1868"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-1.45%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "MCDex"
$ws.Range("C12").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D12").Value = "'11.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'37.95%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09849"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.08%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.04748"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'14.35%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.1064"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.21%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001285"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.48%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04215"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.52%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.005881"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.51%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.371"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-5.54%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3475"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.32%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1421"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'4.41%"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2509"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-3.05%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "'0.001254"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.35%"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "'0.004319"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.77%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001193"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-3.07%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-5.83%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02589"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-2.98%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05504"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.39%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007563"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.58%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1401"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.34%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007471"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-34.69%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002024"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.93%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008362"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-13.25%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00007089"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.44%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.32%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'1.57%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003520"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-1.09%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.32%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002003"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.32%"
$ws.Range("E50").Style = "Normal"

Write-Host "Applied 103 cell changes"